$wb = $excel.ActiveWorkbook

# --- Sheet: ev_charging_uc ---
# Update the two "~TimeSlices" combination strings (C13 and C14).
# Formulas in G7 (=C14) and G8 (=C13) will recalculate automatically.
$wsUc = $wb.Worksheets.Item("ev_charging_uc")
$wsUc.Range("C13").Value = "FaP,SaP,WaD,RaD,SaD,RaP,WaP,FaD"
$wsUc.Range("C14").Value = "SaN,WaN,FaN,FaP,SaP,RaN,WaP,RaP"

# --- Sheet: re_profiles ---
$wsRe = $wb.Worksheets.Item("re_profiles")

# Update N11:N22 values (previously all 0)
$wsRe.Range("N11").Value = 0.097172680668432682
$wsRe.Range("N12").Value = 0.10539924156265808
$wsRe.Range("N13").Value = 0.018540853079282008
$wsRe.Range("N14").Value = 0.13990776695821097
$wsRe.Range("N15").Value = 0.14899437241819788
$wsRe.Range("N16").Value = 0.027122130054051344
$wsRe.Range("N17").Value = 0.068511894999283909
$wsRe.Range("N18").Value = 0.073165268177176404
$wsRe.Range("N19").Value = 0.01532870077005929
$wsRe.Range("N20").Value = 0.14132794813818267
$wsRe.Range("N21").Value = 0.13886151961176466
$wsRe.Range("N22").Value = 0.025667623562496359

# Swap Q13/R13 with Q14/R14 (season "R" and "F" rows swap hydro values)
$wsRe.Range("Q13").Value = "F"
$wsRe.Range("R13").Value = 0.26702915316982878
$wsRe.Range("Q14").Value = "R"
$wsRe.Range("R14").Value = 0.30301943544655252
